# add summary of sprint 1
#
# The "Rejestr Produktu" (Product backlog) sheet already lists the three
# backlog items reserved for "Sprint 2" (rows 14-16, under the "Sprint 2"
# header in row 12). This change adds the Sprint-1 summary numbers next to
# them: effort ("Pracochlonnosc", column C) and business value ("Business
# Value (BV)", column D) for each of those three items, plus a small effort
# total for the first table's last row (D10), and a stray space marker in J7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rejestr Produktu")
$ws.Activate()

# Extra (otherwise-empty) marker cell picked up to the right of the first
# table, on the row of the last fully-populated backlog item.
$ws.Range("J7").Value = " "

# Effort total added next to the "Do pozniejszej negocjacji" row of the
# first table.
$ws.Range("D10").Value = 2

# Sprint-1 summary: effort (C) / business value (D) for each of the three
# Sprint 2 backlog items listed below the first table.
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1

# Move the view/selection down to the newly-filled-in rows, as in the
# authored edit (topLeftCell A5, active cell A14).
[void]$ws.Range("A14").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 5
$aw.ScrollColumn = 1
